$d = $word.ActiveDocument


$rng = $d.Content
$rng.Find.Execute("Clinical Interpretation", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if (-not $rng.Find.Found) { throw "Anchor not found: Clinical Interpretation" }
$p = $rng.Paragraphs(1).Range
$xml_473CF755 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="473CF755" w14:textId="6EF429F6" w:rsidR="00AE3A80" w:rsidRDefault="00851FCF" w:rsidP="006D4A94"><w:pPr><w:spacing w:before="120" w:after="120"/><w:ind w:left="2018" w:right="-8" w:hanging="2023"/><w:jc w:val="both"/><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">Clinical </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Interpretation</w:t></w:r><w:r w:rsidR="00E03BDE"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r><w:r w:rsidR="00E03BDE"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:tab/></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00C61D64"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>CLINICAL_INTERPRETATION1_IN</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p.InsertXML($xml_473CF755)


$rng = $d.Content
$rng.Find.Execute("Please note, variant origin", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if (-not $rng.Find.Found) { throw "Anchor not found: Please note, variant origin" }
$p = $rng.Paragraphs(1).Range
$xml_61B94B9D = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="61B94B9D" w14:textId="7D5E52AF" w:rsidR="009C371C" w:rsidRPr="00B84D96" w:rsidRDefault="00851FCF" w:rsidP="009904B4"><w:pPr><w:keepNext/><w:spacing w:before="120" w:after="120"/><w:ind w:left="2126" w:hanging="2126"/><w:jc w:val="both"/><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Reportable Variants</w:t></w:r><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00490389" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>Please note, variant origin (somatic or germline) cannot be determined by this assay. Variant origin is assumed here based on ancillary information (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>e.g.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> population databases, literature, variant read frequency) for the purpose of clinical interpretation however testing of a germline sample may be recommended in some circumstances.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p.InsertXML($xml_61B94B9D)


$rng = $d.Content
$rng.Find.Execute("AllHaem", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if (-not $rng.Find.Found) { throw "Anchor not found: AllHaem" }
$p = $rng.Paragraphs(1).Range
$xml_2B6B2AFD = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="2B6B2AFD" w14:textId="0EF08D4A" w:rsidR="00B60F7D" w:rsidRPr="00B84D96" w:rsidRDefault="00B60F7D" w:rsidP="009904B4"><w:pPr><w:tabs><w:tab w:val="left" w:pos="8647"/><w:tab w:val="left" w:pos="9540"/></w:tabs><w:spacing w:before="120" w:after="120"/><w:jc w:val="both"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">DNA is analysed by targeted gene sequencing of coding regions and flanking splice sites (within 2 bp) of the genes listed below. Libraries are prepared using a custom QIAGEN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>QIAseq</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> single primer extension-based panel (Peter MacCallum Cancer Centre </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>AllHaem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> v3</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">) and sequenced on an Illumina NextSeq500 with 150 bp paired end reads. </w:t></w:r><w:r w:rsidR="00BC2B51" w:rsidRPr="00BC2B51"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>A customised CLC bioinformatics pipeline including QIAGEN CLC</w:t></w:r><w:r w:rsidR="00BC2B51"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> e</w:t></w:r><w:r w:rsidR="00BC2B51" w:rsidRPr="00BC2B51"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">nterprise </w:t></w:r><w:r w:rsidR="00BC2B51"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>s</w:t></w:r><w:r w:rsidR="00BC2B51" w:rsidRPr="00BC2B51"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">olutions is used to generate aligned reads and call variants (single nucleotide variants and short insertions or deletions) against the hg19 human reference genome. </w:t></w:r><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">Variants are analysed using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>PathOS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> software (Peter Mac) and described according to HGVS nomenclature version 19.01 (http://varnomen.hgvs.org/) with minor differences in accordance with Peter MacCallum Cancer Centre Molecular Pathology departmental policy. The following population variation and cancer or genetic disease databases are commonly used in addition to literature review to assist with variant interpretation: the Genome Aggregation Database (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>gnomAD</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">; gnomad.broadinstitute.org), the Catalogue of Somatic Mutations in Cancer (COSMIC; cancer.sanger.ac.uk), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>ClinVar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> (ncbi.nlm.nih.gov/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>clinvar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>) and the IARC TP53 Database (p53.iarc.fr). V</w:t></w:r><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>ariant origin (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>i.e.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> somatic or germline) is assumed based on ancillary information (e.g. population databases, literature, variant read frequency) for the purpose of clinical interpretation. All assumed somatic variants are reported (and generally considered clinically significant). Variants of uncertain origin are also reported, as are likely benign germline polymorphisms if sufficiently rare and otherwise undescribed. Testing of a non-haematological specimen may be recommended to evaluate variant origin. Recurrent population variants are not </w:t></w:r><w:r w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>reported.</w:t></w:r><w:r w:rsidR="00260F74" w:rsidRPr="00260F74"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> When performed, FLT3-ITDs are tested for by fragment length analysis using capillary electrophoresis. The FLT3 allelic ratio is calculated by peak height ITD/peak height WT.</w:t></w:r><w:r w:rsidR="00260F74" w:rsidRPr="00260F74"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00260F74" w:rsidRPr="00260F74"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>The detection limit for FLT3-ITDs is approximately 1%.</w:t></w:r><w:r w:rsidR="00260F74"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>Somatic variant categorisation</w:t></w:r><w:r w:rsidR="00107A35" w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> (modified from AMP/ASCO/CAP guidelines</w:t></w:r><w:r w:rsidR="002A7213" w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="16"/><w:szCs w:val="18"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>1</w:t></w:r><w:r w:rsidR="00107A35" w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>) –</w:t></w:r><w:r w:rsidR="00107A35" w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>Variants</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> are curated and categorised according to the clinical context of the patient and categorised </w:t></w:r><w:r w:rsidR="00403AE4" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>as</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>DIAGNOSTIC</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> (the variant either defines a diagnostic category or is sufficiently specific for the clinical context to contribute to diagnostic </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>subcategorisation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">), </w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>PROGNOSTIC</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> (the variant has been associated in large trials/series with inferior or superior outcomes in either the context of a specific therapy or independent of therapy</w:t></w:r><w:r w:rsidR="00A44174" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">. Note this does not </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>take into account</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> interaction between prognostic variants present in the individual patient. Relevant pairwise interactions are presented in the</w:t></w:r><w:r w:rsidR="00FE4ED4" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> clinical summary</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">), </w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>DRUG TARGET</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> (the variant or variant class is specifically targeted by a therapeutic agent</w:t></w:r><w:r w:rsidR="00403AE4" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">this category only includes therapeutic agents that are clinically advanced and generally available through either reimbursement or clinical trials </w:t></w:r><w:r w:rsidR="00403AE4" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>[</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>i.e. not early stage investigational agents</w:t></w:r><w:r w:rsidR="00403AE4" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>]</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">), </w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>DRUG RESISTANCE</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> (the variant is specifically associated with resistance to a targeted agent [i.e. does not include non-specific resistance to non-targeted therapies]), </w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>MRD MARKER</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> (the variant is an established </w:t></w:r><w:r w:rsidR="00FE4ED4" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>bio</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">marker </w:t></w:r><w:r w:rsidR="00FE4ED4" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">for which assessment </w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>at MRD sensitivity</w:t></w:r><w:r w:rsidR="00FE4ED4" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> after therapy</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> is accepted practice). If the variant is not categorised into any of the above categories it is assigned </w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B00188"><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>CLONAL MARKER</w:t></w:r><w:r w:rsidR="004D31A3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> indicating its utility in defining the presence of a clonal haematopoietic process</w:t></w:r><w:r w:rsidR="00FE4ED4" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> in the specimen</w:t></w:r><w:r w:rsidR="00935043" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="00201980" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> These </w:t></w:r><w:r w:rsidR="00134296" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>categorisations</w:t></w:r><w:r w:rsidR="00201980" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> are general in nature and</w:t></w:r><w:r w:rsidR="001800BE" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> may not be applicable </w:t></w:r><w:r w:rsidR="00B942E3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">to the </w:t></w:r><w:r w:rsidR="00695366" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">specific </w:t></w:r><w:r w:rsidR="00B942E3" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>clin</w:t></w:r><w:r w:rsidR="00B942E3" w:rsidRPr="00260F74"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>ic</w:t></w:r><w:r w:rsidR="003D798E" w:rsidRPr="00260F74"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>opathological</w:t></w:r><w:r w:rsidR="00B942E3" w:rsidRPr="00260F74"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> context of the patient</w:t></w:r><w:r w:rsidR="00201980" w:rsidRPr="00260F74"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p.InsertXML($xml_2B6B2AFD)


$rng = $d.Content
$rng.Find.Execute("with the exception of CEBPA", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if (-not $rng.Find.Found) { throw "Anchor not found: with the exception of CEBPA" }
$p = $rng.Paragraphs(1).Range
$xml_30AB03C3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="30AB03C3" w14:textId="3D94D113" w:rsidR="0099774F" w:rsidRPr="00B84D96" w:rsidRDefault="00B60F7D" w:rsidP="009904B4"><w:pPr><w:tabs><w:tab w:val="left" w:pos="8647"/><w:tab w:val="left" w:pos="9540"/></w:tabs><w:spacing w:before="120" w:after="120"/><w:jc w:val="both"/><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>The detection limit of this assay for specimens sequenced to the target read depth of 500x is a variant allele frequency (VAF) of approximately 2%</w:t></w:r><w:r w:rsidR="00935043" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> with the exception of CEBPA (detection limit ~ 10%) and ASXL1 c.1934</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>dup;p</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>.Gly646Trpfs*12 (detection limit ~ 5%)</w:t></w:r><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">. This assay is primarily </w:t></w:r><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:lastRenderedPageBreak/><w:t>qualitative however, the variant read frequency (VRF) is provided to assist with variant interpretation and is assumed to approximate VAF in most instances (noting that the VAF of some insertions/deletions may be underrepresented due to assay-based allele bias). The measurement of uncertainty</w:t></w:r><w:r w:rsidR="00620195" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> provided as </w:t></w:r><w:r w:rsidR="007F6B2B" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">a </w:t></w:r><w:r w:rsidR="00620195" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>percentage relative standard uncertainty (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>i.e.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> CV%)</w:t></w:r><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> for variants with VAF</w:t></w:r><w:r w:rsidR="00620195" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>s</w:t></w:r><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> of </w:t></w:r><w:r w:rsidR="00620195" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">5%, 10%-20%, </w:t></w:r><w:r w:rsidR="0099774F" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>3</w:t></w:r><w:r w:rsidR="00620195" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>0</w:t></w:r><w:r w:rsidR="0099774F" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>%</w:t></w:r><w:r w:rsidR="00620195" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>-</w:t></w:r><w:r w:rsidR="0099774F" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>4</w:t></w:r><w:r w:rsidR="00620195" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>0% and 50% are on average, 10.2%, 10.4%, 3.5% and 4.4%, respectively</w:t></w:r><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">. Copy number variations, loss of heterozygosity, structural rearrangements or aneuploidies are not reported. </w:t></w:r><w:r w:rsidR="007F6B2B" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Insertions or deletions (particularly those &gt; 25 bp in length</w:t></w:r><w:r w:rsidR="00077ED4"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>)</w:t></w:r><w:r w:rsidR="007F6B2B" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>, including FLT3-ITDs, are not reliably detected by this assay</w:t></w:r><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">. Genes are analysed using the reference transcripts listed below; coding exons found in alternative transcripts are not assessed by this assay. This assay does not distinguish between somatic and germline variants. </w:t></w:r><w:r w:rsidR="0064555A"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>In addition, the clonal origin of somatic variants (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>i.e.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> disease compartment or cell lineage) cannot be determined. </w:t></w:r><w:r w:rsidR="004C2468" w:rsidRPr="004C2468"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>S</w:t></w:r><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">ynonymous variants </w:t></w:r><w:r w:rsidR="00D22DE6" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>are</w:t></w:r><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> not</w:t></w:r><w:r w:rsidR="00D22DE6" w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> routinely</w:t></w:r><w:r w:rsidRPr="00B84D96"><w:rPr><w:rFonts w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> reported. Please note Peter Mac assumes sample identification, family relationships, and clinical diagnoses are as stated on the request. Our clinical recommendations may be based on evidence from third-party data sources and should be interpreted in the context of all other clinical and laboratory information for this patient.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p.InsertXML($xml_30AB03C3)


$rng = $d.Content
$rng.Find.Execute("Please note FLT3-ITDs are not detected", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if (-not $rng.Find.Found) { throw "Anchor not found: Please note FLT3-ITDs are not detected" }
$p = $rng.Paragraphs(1).Range
$xml_4CE5A378 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="4CE5A378" w14:textId="2BAE8691" w:rsidR="00260F74" w:rsidRPr="00260F74" w:rsidRDefault="00045DDC" w:rsidP="00260F74"><w:pPr><w:tabs><w:tab w:val="left" w:pos="8647"/><w:tab w:val="left" w:pos="9540"/></w:tabs><w:spacing w:after="120"/><w:jc w:val="both"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:szCs w:val="12"/></w:rPr><w:t>*</w:t></w:r><w:r w:rsidR="00260F74" w:rsidRPr="00260F74"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> Please note FLT3-ITDs are not detected with this assay. A separate assay may have been </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:szCs w:val="16"/></w:rPr><w:t>performed,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> result included in Test Results if sample tested. </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:szCs w:val="16"/></w:rPr><w:t>^</w:t></w:r><w:r w:rsidR="00260F74" w:rsidRPr="00260F74"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:szCs w:val="12"/><w:vertAlign w:val="superscript"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00260F74" w:rsidRPr="00260F74"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:szCs w:val="16"/></w:rPr><w:t>Partial coverage of region</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p.InsertXML($xml_4CE5A378)


$d.Content.Find.Execute("1-Nov-2023", $true, $false, $false, $false, $false, $true, 1, $false, "15-Nov-2023", 2) | Out-Null
